$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated coin prices / 1h volume percentages (and, for rows 39-40 and
# 48/50/51, the coin name + link that moved to a new rank position) to
# match the latest coinranking.com snapshot.
#
# Price-column values that would otherwise be auto-parsed as numbers by
# Excel are written with the cell briefly marked as Text (NumberFormat
# "@") so they keep their original string form (e.g. "0.200", "41.76")
# exactly as the source data stores them; the style is then reverted so
# no visible formatting change is introduced.

# Row 2
$ws.Range("D2").Value = '94.284.73'
$ws.Range("E2").Value = '  -3.97%  '
# Row 3
$ws.Range("D3").Value = '3.415.92'
$ws.Range("E3").Value = '  +0.75%  '
# Row 4
$ws.Range("E4").Value = '  +0.10%  '
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '236.95'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -6.70%  '
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '641.77'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.25%  '
# Row 7
$ws.Range("E7").Value = '  -0.71%  '
# Row 8
$ws.Range("E8").Value = '  -4.84%  '
# Row 9
$ws.Range("E9").Value = '  +0.14%  '
# Row 10
$ws.Range("E10").Value = '  -7.13%  '
# Row 11
$ws.Range("D11").Value = '3.411.90'
$ws.Range("E11").Value = '  +0.69%  '
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.200'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.37%  '
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '41.76'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.34%  '
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.21'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.73%  '
# Row 15
$ws.Range("D15").Value = '94.073.76'
$ws.Range("E15").Value = '  -3.83%  '
# Row 16
$ws.Range("D16").Value = '4.061.50'
$ws.Range("E16").Value = '  +1.17%  '
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000251'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.72%  '
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '8.33'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.46%  '
# Row 19
$ws.Range("D19").Value = '3.419.66'
$ws.Range("E19").Value = '  +0.97%  '
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.13%  '
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '11.60'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.86%  '
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.494'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -6.57%  '
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '498.55'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.64%  '
# Row 24
$ws.Range("E24").Value = '  -6.41%  '
# Row 25
$ws.Range("E25").Value = '  -4.52%  '
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.51'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -6.97%  '
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '91.29'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.58%  '
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '11.93'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.50%  '
# Row 29
$ws.Range("D29").Value = '3.603.68'
$ws.Range("E29").Value = '  +1.01%  '
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '11.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.16%  '
# Row 31
$ws.Range("E31").Value = '  +0.17%  '
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.76'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +6.48%  '
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.138'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -3.36%  '
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.998'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.23%  '
# Row 35
$ws.Range("E35").Value = '  -4.34%  '
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '29.53'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.04%  '
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.550'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.92%  '
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '552.81'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.38%  '
# Row 39
$ws.Range("B39").Value = 'Fetch.AI'
$ws.Range("C39").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.45'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.43%  '
# Row 40
$ws.Range("B40").Value = 'RenderToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '7.63'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.64%  '
# Row 41
$ws.Range("E41").Value = '  +0.03%  '
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.150'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.24%  '
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.901'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.20%  '
# Row 44
$ws.Range("E44").Value = '  -1.45%  '
# Row 45
$ws.Range("E45").Value = '  +0.24%  '
# Row 46
$ws.Range("E46").Value = '  -1.33%  '
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.59'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.56%  '
# Row 48
$ws.Range("B48").Value = 'OKB'
$ws.Range("C48").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.48'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.07%  '
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0411'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.22%  '
# Row 50
$ws.Range("B50").Value = 'dogwifhat'
$ws.Range("C50").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.32'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +3.22%  '
# Row 51
$ws.Range("B51").Value = 'Stacks'
$ws.Range("C51").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.18'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.33%  '
